# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.143.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "'1.825.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'310.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "'0.4959"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.67%  "
$ws.Range("D8").Value = "'0.3922"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("D9").Value = "'0.09828"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +24.94%  "
$ws.Range("D10").Value = "'1.110"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").Value = "'41.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "'6.457"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "'1.001"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").Value = "'1.820.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "'7.322"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").Value = "  +5.56%  "
$ws.Range("D18").Value = "'92.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").Value = "'0.06663"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").Value = "'6.019"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "'28.191.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("D24").Value = "'11.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").Value = "'2.247"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("D26").Value = "'158.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'20.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("D28").Value = "'2.033.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("D29").Value = "'2.429"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").Value = "'127.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.68%  "
$ws.Range("D31").Value = "'0.1054"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("D32").Value = "'1.041"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").Value = "'5.610"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D35").Value = "'0.06739"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.34%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'9.025"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02347"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("D39").Value = "'4.974"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  -1.94%  "
$ws.Range("D41").Value = "'0.6227"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "'1.184"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").Value = "'13.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("E45").Value = "  -1.62%  "
$ws.Range("D46").Value = "'3.710"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").Value = "'1.276"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("D48").Value = "'124.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").Value = "'1.953"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").Value = "'1.184"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.98%  "
$ws.Range("D51").Value = "'0.06791"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.98%  "
